$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Score" label + its formula move up from old row 5 into row 4 (C4/D4)
$ws.Range("C4").Value = "Score"
$ws.Range("D4").Formula = '=(D2-D3)*$B$4+$B$5'

# Row 5: new "Probability" label + formula (previously lived in row 6, now referencing D4)
$ws.Range("C5").Value = "Probability"
$ws.Range("D5").Formula = '=$B$7*EXP(D4)/(1+$B$7*EXP(D4))'

# Row 6: new formula that thresholds the probability
$ws.Range("D6").Formula = '=IF(D5<0.5, "Can not be determined", D5)'

# Update the last-used selection
$ws.Range("E10").Select()
